$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (shifts old row 4 -> row 5)
$ws.Rows(4).Insert()

# --- Row 2 ---
$ws.Range("A2").Value = "Bulgaria"
$ws.Range("B2").Value = "'3"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = "Bank (Money Center)"
$ws.Range("D2").Value = 0.0747
$ws.Range("E2").Value = 0.182
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 74.14399999999999
$ws.Range("L2").Value = 0.2335979836168872
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("U2").Value = 1107.5
$ws.Range("V2").Value = 3.659947124917382
$ws.Range("W2").Value = 0.02111781517942204
$ws.Range("X2").Value = 0.05616533060496844
$ws.Range("Y2").Value = -0.03504751542554639
$ws.Range("Z2").Value = -8.82058692752333
$ws.Range("AA2").Value = -0
$ws.Range("AB2").Value = 0.04611804328932199
$ws.Range("AC2").Value = -0.04611804328932199
$ws.Range("AD2").Value = 366.2
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 366.2
$ws.Range("AG2").Value = -741.3
$ws.Range("AH2").Value = 0.5475478468899522
$ws.Range("AI2").Value = 0.2490478781284004
$ws.Range("AJ2").Value = 1.689765215409164
$ws.Range("AK2").Value = -2.042711490768806
$ws.Range("AL2").Value = 0
$ws.Range("AM2").Value = 0
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()

# --- Row 3 ---
$ws.Range("A3").Value = "Bulgaria"
$ws.Range("B3").Value = "Texim Bank AD (BUL:5CP)"
$ws.Range("C3").Value = "Bank (Money Center)"
$ws.Range("D3").Value = 0.09810000000000001
$ws.Range("E3").Value = 0.182
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0.194
$ws.Range("L3").Value = 0.02042105263157895
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = -0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("S3").Value = 0
$ws.Range("U3").Value = 49
$ws.Range("V3").Value = 0.9210526315789473
$ws.Range("W3").Value = 0.009194312796208531
$ws.Range("X3").Value = 0.04609342934023551
$ws.Range("Y3").Value = -0.03689911654402698
$ws.Range("Z3").Value = -21.44469525959377
$ws.Range("AA3").Value = -0
$ws.Range("AB3").Value = 0.0428707905552534
$ws.Range("AC3").Value = -0.0428707905552534
$ws.Range("AD3").Value = 13.3
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 13.3
$ws.Range("AG3").Value = -35.7
$ws.Range("AH3").Value = 0.2
$ws.Range("AI3").Value = 0.3653846153846154
$ws.Range("AJ3").Value = -2.04
$ws.Range("AK3").Value = 2.833333333333333
$ws.Range("AL3").Value = 0
$ws.Range("AM3").Value = 0
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()

# --- Row 4 ---
$ws.Range("A4").Value = "Bulgaria"
$ws.Range("B4").Value = "Central Cooperative Bank AD (BUL:4CF)"
$ws.Range("C4").Value = "Bank (Money Center)"
$ws.Range("D4").Value = 0.029
$ws.Range("E4").Value = 0.171
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 6.65
$ws.Range("L4").Value = 0.0689119170984456
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = -0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = -0
$ws.Range("S4").Value = 0
$ws.Range("U4").Value = 957.3
$ws.Range("V4").Value = 12.28883183568678
$ws.Range("W4").Value = 0.02111781517942204
$ws.Range("X4").Value = 0.05616533060496844
$ws.Range("Y4").Value = -0.03504751542554639
$ws.Range("Z4").Value = -0.1583680703712213
$ws.Range("AA4").Value = -0
$ws.Range("AB4").Value = 0.04611804328932199
$ws.Range("AC4").Value = -0.04611804328932199
$ws.Range("AD4").Value = 48.5
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 48.5
$ws.Range("AG4").Value = -908.8
$ws.Range("AH4").Value = 0.3837025316455696
$ws.Range("AI4").Value = 0.1214625594790884
$ws.Range("AJ4").Value = 1.093753760982068
$ws.Range("AK4").Value = 1.628673835125448
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = 0

# --- Row 5 ---
$ws.Range("A5").Value = "Bulgaria"
$ws.Range("B5").Value = "First Investment Bank AD (BUL:5F4)"
$ws.Range("C5").Value = "Bank (Money Center)"
$ws.Range("D5").Value = 0.0747
$ws.Range("E5").Value = 0.416
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 67.3
$ws.Range("L5").Value = 0.3183538315988647
$ws.Range("M5").Value = -0
$ws.Range("N5").Value = -0
$ws.Range("O5").Value = -0
$ws.Range("P5").Value = -0
$ws.Range("Q5").Value = -0
$ws.Range("R5").Value = -0
$ws.Range("S5").Value = 0
$ws.Range("U5").Value = 101.2
$ws.Range("V5").Value = 0.5900874635568514
$ws.Range("W5").Value = 0.13188320595728
$ws.Range("X5").Value = 0.08731511767192865
$ws.Range("Y5").Value = 0.04456808828535139
$ws.Range("Z5").Value = 0.3684216947049402
$ws.Range("AA5").Value = 0
$ws.Range("AB5").Value = 0.0506419969345675
$ws.Range("AC5").Value = -0.0506419969345675
$ws.Range("AD5").Value = 304.4
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 304.4
$ws.Range("AG5").Value = 203.2
$ws.Range("AH5").Value = 0.6396301744063879
$ws.Range("AI5").Value = 0.2941915531071809
$ws.Range("AJ5").Value = 0.5423005070723245
$ws.Range("AK5").Value = 0.2176754151044456
$ws.Range("AL5").Value = 0
$ws.Range("AM5").Value = 0

Write-Host "Update complete"
